$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.488.64'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +1.57%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.674.22'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '220.52'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.29%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5309'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.72%  '
$ws.Range('E7').Value = '  +0.00%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2689'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +2.89%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06385'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E11').Value = '  +1.78%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.711.82'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +3.58%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.492'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.52%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5574'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.44%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0₅8335'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.91%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.61'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.92%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '26.505.39'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.59%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.764'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.63%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '192.54'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.32%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.34'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.16%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.312'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.27%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.1275'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +4.57%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '138.46'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -5.27%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.408'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  +2.82%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.429'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.46%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.06280'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +5.54%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.290'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +1.89%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.605'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +5.88%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.421'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.65%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.691'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.53%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.009'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.48%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.6151'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +9.13%  '
$ws.Range('E36').Value = '  +1.08%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.783'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '6.122'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +4.85%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01620'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.092.19'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +6.10%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.8614'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  -0.13%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '100.71'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.01%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.820.30'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.39%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '58.79'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +5.19%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '8.191'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('E47').Value = '  +0.56%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.523'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +9.90%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.05197'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.64%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '6.010'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('E51').Value = '  +0.29%  '
